{"js": "// Quality Management Plan \u2013 \"Quality Control Measurements\" section:\n//  1) Narrow the right indent of every paragraph in the section\n//     (w:ind right 873 -> 4 twips, i.e. ~43.65pt -> 0.2pt).\n//  2) In the closing \"In summary, ...\" paragraph, the rendered page\n//     break now falls right before \"Continuous assessment ...\" instead\n//     of before \"In summary, ...\", so the run is split there and the\n//     <w:lastRenderedPageBreak/> marker moves to the new run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the start of the \"Quality Control Measurements\" narrative block\n// by its distinctive opening sentence, then walk the fixed-size run of\n// paragraphs that make up that block (5 text paragraphs + 4 blank\n// spacer paragraphs interleaved between them).\nconst anchorText =\n  \"will leverage Agile and Scrum methodologies to foster continuous inspection\";\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    startIndex = i;\n    break;\n  }\n}\nif (startIndex === -1) {\n  throw new Error(\"Could not locate the Quality Control Measurements block\");\n}\n\nconst blockSize = 9;\nconst blockParagraphs = [];\nfor (let i = startIndex; i < startIndex + blockSize; i++) {\n  blockParagraphs.push(paragraphs.items[i]);\n}\n\n// 1) Tighten the right indent on every paragraph in the block.\nfor (const p of blockParagraphs) {\n  p.rightIndent = 0.2; // points (4 twips)\n}\n\n// 2) Move the page-break hint in the final paragraph from before\n//    \"In summary, the\" to before \"Continuous assessment ...\", by\n//    replacing the paragraph's OOXML with the run split applied.\nconst lastParagraph = blockParagraphs[blockParagraphs.length - 1];\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  \"</Relationships></pkg:xmlData></pkg:part>\" +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  \"<w:body>\" +\n  '<w:p w14:paraId=\"41C8F396\" w14:textId=\"52B97277\" w:rsidR=\"004D2246\" w:rsidRDefault=\"004C5926\" w:rsidP=\"004C5926\">' +\n  '<w:pPr><w:ind w:left=\"1440\" w:right=\"4\" w:hanging=\"450\"/></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">In summary, the </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>SurveiRams</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> Ticketing System project will adopt Agile and Scrum methodologies to </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/><w:r><w:t>establish</w:t></w:r><w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> a collaborative and dynamic quality control strategy. </w:t></w:r>' +\n  \"<w:r><w:lastRenderedPageBreak/><w:t>Continuous assessment of the product's quality will be performed, with regular improvements implemented. All quality control measurements will be collected and tracked on a shared platform in real-time. The team will collaborate to address any issues and drive necessary enhancements.</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nlastParagraph.insertOoxml(newParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Quality Management Plan - \"Quality Control Measurements\" section:\n#  1) Narrow the right indent of every paragraph in the section\n#     (w:ind right 873 -> 4 twips, i.e. ~43.65pt -> 0.2pt).\n#  2) In the closing \"In summary, ...\" paragraph, the rendered page\n#     break now falls right before \"Continuous assessment ...\" instead\n#     of before \"In summary, ...\", so the run is split there and the\n#     lastRenderedPageBreak marker moves to the new run.\n\n$d = $word.ActiveDocument\n\n# Locate the start of the \"Quality Control Measurements\" narrative block\n# by its distinctive opening sentence, then walk the fixed-size run of\n# paragraphs that make up that block (5 text paragraphs + 4 blank\n# spacer paragraphs interleaved between them = 9 paragraphs total).\n$startIdx = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n  $i = $i + 1\n  if ($p.Range.Text -like \"*will leverage Agile and Scrum methodologies to foster continuous inspection*\") {\n    $startIdx = $i\n    break\n  }\n}\nif ($startIdx -eq -1) {\n  throw \"Could not locate the Quality Control Measurements block\"\n}\n\n$blockSize = 9\n$endIdx = $startIdx + $blockSize - 1\n\n# 1) Tighten the right indent on every paragraph in the block, and\n#    remember the last paragraph of the block for step 2.\n$lastPara = $null\n$j = 0\nforeach ($p in $d.Paragraphs) {\n  $j = $j + 1\n  if ($j -ge $startIdx -and $j -le $endIdx) {\n    $p.Range.ParagraphFormat.RightIndent = 0.2\n    $lastPara = $p\n  }\n  if ($j -gt $endIdx) {\n    break\n  }\n}\n\n# 2) Move the page-break hint in the final paragraph from before\n#    \"In summary, the\" to before \"Continuous assessment ...\", by\n#    replacing the paragraph's content with the run split applied.\n$newParagraphXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" w14:paraId=\"41C8F396\" w14:textId=\"52B97277\" w:rsidR=\"004D2246\" w:rsidRDefault=\"004C5926\" w:rsidP=\"004C5926\"><w:pPr><w:ind w:left=\"1440\" w:right=\"4\" w:hanging=\"450\"/></w:pPr><w:r><w:t xml:space=\"preserve\">In summary, the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>SurveiRams</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> Ticketing System project will adopt Agile and Scrum methodologies to </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>establish</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> a collaborative and dynamic quality control strategy. </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>Continuous assessment of the product''s quality will be performed, with regular improvements implemented. All quality control measurements will be collected and tracked on a shared platform in real-time. The team will collaborate to address any issues and drive necessary enhancements.</w:t></w:r></w:p>'\n\n$null = $lastPara.Range.InsertXML($newParagraphXml)\n\nWrite-Output \"Quality Control Measurements section updated\"\n"}
